# Clemson Dashboard update: EOW Oct 31 up to Nov 2
$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet 1: "Daily Data"
# ------------------------------------------------------------------
$daily = $wb.Worksheets.Item("Daily Data")

# Row 70 (2020-11-01 / serial 44134, date already present in A70)
$daily.Range("B70").Value = 21
$daily.Range("C70").FormulaR1C1 = "=AVERAGE(R[-6]C[-1]:RC[-1])"
$daily.Range("D70").FormulaR1C1 = "=(RC[-2]/25822)*100000"
$daily.Range("E70").Value = 2189
$daily.Range("F70").FormulaR1C1 = "=(RC[-1]/25822)*100000"
$daily.Range("G70").FormulaR1C1 = "=RC[-5]/RC[-2]"
$daily.Range("H70").FormulaR1C1 = "=AVERAGE(R[-6]C[-1]:RC[-1])"
$daily.Range("I70").FormulaR1C1 = "=RC[-7]+R[-1]C"
$daily.Range("J70").FormulaR1C1 = "=R[-1]C+RC[-5]"
$daily.Range("S70").Formula = "=IF(K70+M70=B70,""EQUAL"",""DIFFER"")"
$daily.Range("T70").Formula = "=IF(L70+N70=E70,""EQUAL"",""DIFFER"")"
$daily.Range("U70").Formula = "=IF(O70+Q70=B70,""EQUAL"",""DIFFER"")"
$daily.Range("V70").Formula = "=IF(P70+R70=E70,""EQUAL"",""DIFFER"")"

# Row 71 (2020-11-02 / serial 44135, date already present in A71)
$daily.Range("B71").Value = 1
$daily.Range("C71").FormulaR1C1 = "=AVERAGE(R[-6]C[-1]:RC[-1])"
$daily.Range("D71").FormulaR1C1 = "=(RC[-2]/25822)*100000"
$daily.Range("E71").Value = 2
$daily.Range("F71").FormulaR1C1 = "=(RC[-1]/25822)*100000"
$daily.Range("G71").FormulaR1C1 = "=RC[-5]/RC[-2]"
$daily.Range("H71").FormulaR1C1 = "=AVERAGE(R[-6]C[-1]:RC[-1])"
$daily.Range("I71").FormulaR1C1 = "=RC[-7]+R[-1]C"
$daily.Range("J71").FormulaR1C1 = "=R[-1]C+RC[-5]"
$daily.Range("S71").Formula = "=IF(K71+M71=B71,""EQUAL"",""DIFFER"")"
$daily.Range("T71").Formula = "=IF(L71+N71=E71,""EQUAL"",""DIFFER"")"
$daily.Range("U71").Formula = "=IF(O71+Q71=B71,""EQUAL"",""DIFFER"")"
$daily.Range("V71").Formula = "=IF(P71+R71=E71,""EQUAL"",""DIFFER"")"

# View: Daily Data tab is selected again, selection moved to F10
$daily.Activate()
$daily.Range("F10").Select()

Write-Host "sheet1 done"
